# Update "想去人数" (F column) figures for both the "展览" sheet and the
# "全部类型" sheet, to reflect freshly scraped attendance counts.
#
# The two sheets list (mostly) the same events but at different row
# offsets (展览 excludes two 演出/performance rows that 全部类型 includes),
# so each sheet is addressed by its own explicit cell references.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Cell -> new value, for the "展览" sheet
$exhibitUpdates = @{
    "F2"  = 1090
    "F3"  = 409
    "F5"  = 140
    "F6"  = 12024
    "F8"  = 80
    "F9"  = 11759
    "F10" = 4743
    "F11" = 532
    "F12" = 70
    "F15" = 926
}

foreach ($cell in $exhibitUpdates.Keys) {
    $wsExhibit.Range($cell).Value = $exhibitUpdates[$cell]
}

# Cell -> new value, for the "全部类型" sheet
$allUpdates = @{
    "F2"  = 1090
    "F3"  = 409
    "F5"  = 140
    "F8"  = 12024
    "F10" = 80
    "F11" = 11759
    "F12" = 4743
    "F13" = 532
    "F14" = 70
    "F17" = 926
}

foreach ($cell in $allUpdates.Keys) {
    $wsAll.Range($cell).Value = $allUpdates[$cell]
}
